{"js": "// Update the worksheet date and the 24 division problems to the new values.\nconst replacements = [\n  [\"2025-08-14 Thursday\", \"2025-08-15 Friday\"],\n  [\"93\u00f79=\", \"72\u00f74=\"],\n  [\"75\u00f78=\", \"15\u00f72=\"],\n  [\"49\u00f72=\", \"95\u00f76=\"],\n  [\"31\u00f72=\", \"12\u00f73=\"],\n  [\"63\u00f76=\", \"84\u00f75=\"],\n  [\"46\u00f79=\", \"40\u00f79=\"],\n  [\"92\u00f78=\", \"12\u00f76=\"],\n  [\"18\u00f78=\", \"13\u00f72=\"],\n  [\"43\u00f76=\", \"58\u00f75=\"],\n  [\"94\u00f74=\", \"44\u00f74=\"],\n  [\"58\u00f76=\", \"24\u00f77=\"],\n  [\"71\u00f78=\", \"62\u00f77=\"],\n  [\"27\u00f72=\", \"27\u00f79=\"],\n  [\"20\u00f78=\", \"30\u00f78=\"],\n  [\"11\u00f77=\", \"15\u00f76=\"],\n  [\"71\u00f74=\", \"59\u00f76=\"],\n  [\"81\u00f75=\", \"12\u00f77=\"],\n  [\"89\u00f73=\", \"84\u00f77=\"],\n  [\"29\u00f77=\", \"52\u00f78=\"],\n  [\"74\u00f76=\", \"28\u00f74=\"],\n  [\"90\u00f77=\", \"57\u00f75=\"],\n  [\"72\u00f79=\", \"67\u00f73=\"],\n  [\"36\u00f77=\", \"76\u00f77=\"],\n  [\"48\u00f72=\", \"31\u00f73=\"],\n  [\"37\u00f79=\", \"61\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 24 division problems to the new values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-08-14 Thursday\", \"2025-08-15 Friday\"),\n    @(\"93\u00f79=\", \"72\u00f74=\"),\n    @(\"75\u00f78=\", \"15\u00f72=\"),\n    @(\"49\u00f72=\", \"95\u00f76=\"),\n    @(\"31\u00f72=\", \"12\u00f73=\"),\n    @(\"63\u00f76=\", \"84\u00f75=\"),\n    @(\"46\u00f79=\", \"40\u00f79=\"),\n    @(\"92\u00f78=\", \"12\u00f76=\"),\n    @(\"18\u00f78=\", \"13\u00f72=\"),\n    @(\"43\u00f76=\", \"58\u00f75=\"),\n    @(\"94\u00f74=\", \"44\u00f74=\"),\n    @(\"58\u00f76=\", \"24\u00f77=\"),\n    @(\"71\u00f78=\", \"62\u00f77=\"),\n    @(\"27\u00f72=\", \"27\u00f79=\"),\n    @(\"20\u00f78=\", \"30\u00f78=\"),\n    @(\"11\u00f77=\", \"15\u00f76=\"),\n    @(\"71\u00f74=\", \"59\u00f76=\"),\n    @(\"81\u00f75=\", \"12\u00f77=\"),\n    @(\"89\u00f73=\", \"84\u00f77=\"),\n    @(\"29\u00f77=\", \"52\u00f78=\"),\n    @(\"74\u00f76=\", \"28\u00f74=\"),\n    @(\"90\u00f77=\", \"57\u00f75=\"),\n    @(\"72\u00f79=\", \"67\u00f73=\"),\n    @(\"36\u00f77=\", \"76\u00f77=\"),\n    @(\"48\u00f72=\", \"31\u00f73=\"),\n    @(\"37\u00f79=\", \"61\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\n$d.Saved = $false\n"}
